# edit.ps1 - applies the "Add styles for payment history and Islamic
# properties summary" change to the finance workbook:
#   * Debts sheet gains paymentHistory / historicalPaid columns and the
#     existing "paid" column is relocated to the end with an updated value.
#   * Two brand-new sheets are appended: "Properties" and "Islamic".
#   * The Expenses sheet gets an updated recurring-info row plus a new
#     (duplicate-ish) expense row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Debts sheet (existing sheet, 12th tab) - re-layout the columns.
# ---------------------------------------------------------------------
$debts = $wb.Worksheets.Item("Debts")

# Header row
$debts.Range("D1").Value = "paymentHistory"
$debts.Range("E1").Value = "historicalPaid"
$debts.Range("F1").Value = "creditor"
$debts.Range("G1").Value = "dueDate"
$debts.Range("H1").Value = "description"
$debts.Range("I1").Value = "paid"

# Data row 2
$debts.Range("D2").Value = '[{"amount":5000,"date":"2026-01-02","note":"","method":"Cash","recordedAt":"2026-01-02T06:58:03.990Z"}]'
$debts.Range("E2").Value = 0
$debts.Range("F2").Value = ""
$debts.Range("G2").NumberFormat = "@"
$debts.Range("G2").Value = "2026-01-02"
$debts.Range("H2").Value = ""
$debts.Range("I2").Value = 5000

# ---------------------------------------------------------------------
# 2. Expenses sheet (2nd tab) - update row 2 recurring info + add row 4.
# ---------------------------------------------------------------------
$expenses = $wb.Worksheets.Item("Expenses")

$expenses.Range("K2").Value = $true
$expenses.Range("L2").Value = "monthly"
$expenses.Range("M2").NumberFormat = "@"
$expenses.Range("M2").Value = "2026-01-31"

$expenses.Range("A4").Value = "office ar chele"
$expenses.Range("B4").Value = "Other"
$expenses.Range("C4").Value = 5000
$expenses.Range("D4").NumberFormat = "@"
$expenses.Range("D4").Value = "2026-01-01"
$expenses.Range("E4").Value = "Hand Cash"
$expenses.Range("F4").Value = "Husna"
$expenses.Range("G4").Value = '["vaiya","me"]'
$expenses.Range("H4").Value = "equal"
$expenses.Range("I4").Value = ""
$expenses.Range("J4").Value = ""
$expenses.Range("K4").Value = $true
$expenses.Range("L4").Value = "monthly"
$expenses.Range("M4").Value = ""

# ---------------------------------------------------------------------
# 3. New "Properties" sheet appended after "Debts".
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$properties = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$properties.Name = "Properties"

$properties.Range("A1").Value = "name"
$properties.Range("B1").Value = "type"
$properties.Range("C1").Value = "value"
$properties.Range("D1").Value = "estimatedValue"
$properties.Range("E1").Value = "location"
$properties.Range("F1").Value = "area"
$properties.Range("G1").Value = "acquiredDate"
$properties.Range("H1").Value = "zakatEligible"
$properties.Range("I1").Value = "notes"

$properties.Range("A2").Value = "gold"
$properties.Range("B2").Value = "Gold"
$properties.Range("C2").Value = 50000
$properties.Range("D2").Value = 50000
$properties.Range("E2").Value = ""
$properties.Range("F2").Value = -0.01
$properties.Range("G2").NumberFormat = "@"
$properties.Range("G2").Value = "2026-01-01"
$properties.Range("H2").Value = $true
$properties.Range("I2").Value = ""

# ---------------------------------------------------------------------
# 4. New, empty "Islamic" sheet appended after "Properties".
# ---------------------------------------------------------------------
$islamic = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $properties)
$islamic.Name = "Islamic"
